$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 11).Value = 1
}
